# Fixbug: import học viên
# Update the report title (drop the "K29 HÀ NỘI + TÂY BẮC" suffix) and
# move the active selection to E5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title in A1 to the generic wording (no more cohort/location suffix)
$ws.Range("A1").Value = "DANH SÁCH HỘI ĐỒNG CHẤM LUẬN VĂN"

# Move / restore the active selection on the sheet
$ws.Range("E5").Select()
